# B1/B2 PowerPoint - Wed, Jul 29, 2020  9:05:50 AM
#
# This edit does two things:
#
#   1. On slide 5, the financial-documents table's style is switched
#      from the deck's custom "Table_0" style
#      ({A36ED6EB-FD7C-4D25-AD11-018755DBB355}) to the built-in table
#      style {582E9031-CDE1-4DA1-B6C6-9BB9AAF292D6}.
#
#   2. The presentation's theme colour scheme (currently the colourful
#      "Integral" / "Red Violet" palette used by the slide master) is
#      recoloured to the plain default "Office Theme" / "Office"
#      palette, i.e. the deck's design swaps from "Integral" back to
#      the built-in "Office Theme" colours.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 --------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{582E9031-CDE1-4DA1-B6C6-9BB9AAF292D6}")
    }
}

# --- 2. Theme colours -------------------------------------------------
function HexToOleRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Colors(1..12):
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = HexToOleRgb($officeThemeColors[$i - 1])
}

# Best-effort: some hosts also let the theme / colour-scheme display
# names be updated to match ("Office Theme" / "Office"). Harmless if
# unsupported.
try { $theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
